$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 101, shifting existing rows 101-125 down to 102-126
$ws.Rows.Item(101).Insert()

# Fill the new row 101 with the new weekly data entry.
# Columns A,B,C,E,F,G,H,I,J,K,L,Q,R,T keep the same values as the (now-shifted) row 102,
# while D, M, N, O, P, S hold the new data.
$ws.Range("A101").Value = 10
$ws.Range("B101").Value = "Vega Modelo de Temuco"
$ws.Range("C101").Value = "La Araucanía"
$ws.Range("D101").Value = 45204
$ws.Range("E101").Value = 9
$ws.Range("F101").Value = "Fruta"
$ws.Range("G101").Value = 100108
$ws.Range("H101").Value = "Tropicales y subtropicales"
$ws.Range("I101").Value = 100108004
$ws.Range("J101").Value = "Papaya"
$ws.Range("K101").Value = "Cultivar IV Región"
$ws.Range("L101").Value = "Primera"
$ws.Range("M101").Value = 120
$ws.Range("N101").Value = 24000
$ws.Range("O101").Value = 24000
$ws.Range("P101").Value = 24000
$ws.Range("Q101").Value = "$/bandeja 10 kilos"
$ws.Range("R101").Value = "Provincia del Elquí"
$ws.Range("S101").Value = 2400
$ws.Range("T101").Value = 10
